# Adds macro-photography scale bar distance rows for a new camera body
# (Canon EOS R5) across three lens/objective groups: Laowa 25mm f/2.8,
# Laowa 85mm f/5.6, and Mitutoyo objectives. Mirrors the existing table's
# structure (one data block per lens, last row of each block using the
# bold/bordered "group total" style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stamp formatting for the new rows by copying it from existing,
#        structurally-identical blocks, then overwrite with new content. ---

# Rows 80-85 (Laowa 25mm f/2.8, 2.5x..5.0x) <- copy from rows 66-71
$ws.Range("A66:I71").Copy() | Out-Null
$ws.Range("A80:I85").PasteSpecial(-4122) | Out-Null

# Rows 86-90 (Laowa 85mm f/5.6, 0.5x..2.0x) <- copy from rows 18-22
$ws.Range("A18:I22").Copy() | Out-Null
$ws.Range("A86:I90").PasteSpecial(-4122) | Out-Null

# Rows 91-92 (Mitutoyo objectives) <- copy from rows 21-22
$ws.Range("A21:I22").Copy() | Out-Null
$ws.Range("A91:I92").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- 2. Fill in the new data (column-major for the repeated camera/image
#        fields, then per-row for the rest) ---
$ws.Range("A80").Value2 = "Canon EOS R5"
$ws.Range("A81").Value2 = "Canon EOS R5"
$ws.Range("A82").Value2 = "Canon EOS R5"
$ws.Range("A83").Value2 = "Canon EOS R5"
$ws.Range("A84").Value2 = "Canon EOS R5"
$ws.Range("A85").Value2 = "Canon EOS R5"
$ws.Range("A86").Value2 = "Canon EOS R5"
$ws.Range("A87").Value2 = "Canon EOS R5"
$ws.Range("A88").Value2 = "Canon EOS R5"
$ws.Range("A89").Value2 = "Canon EOS R5"
$ws.Range("A90").Value2 = "Canon EOS R5"
$ws.Range("A91").Value2 = "Canon EOS R5"
$ws.Range("A92").Value2 = "Canon EOS R5"
$ws.Range("E80").Value2 = "8192"
$ws.Range("E81").Value2 = "8192"
$ws.Range("E82").Value2 = "8192"
$ws.Range("E83").Value2 = "8192"
$ws.Range("E84").Value2 = "8192"
$ws.Range("E85").Value2 = "8192"
$ws.Range("E86").Value2 = "8192"
$ws.Range("E87").Value2 = "8192"
$ws.Range("E88").Value2 = "8192"
$ws.Range("E89").Value2 = "8192"
$ws.Range("E90").Value2 = "8192"
$ws.Range("E91").Value2 = "8192"
$ws.Range("E92").Value2 = "8192"
$ws.Range("F80").Value2 = "5464"
$ws.Range("F81").Value2 = "5464"
$ws.Range("F82").Value2 = "5464"
$ws.Range("F83").Value2 = "5464"
$ws.Range("F84").Value2 = "5464"
$ws.Range("F85").Value2 = "5464"
$ws.Range("F86").Value2 = "5464"
$ws.Range("F87").Value2 = "5464"
$ws.Range("F88").Value2 = "5464"
$ws.Range("F89").Value2 = "5464"
$ws.Range("F90").Value2 = "5464"
$ws.Range("F91").Value2 = "5464"
$ws.Range("F92").Value2 = "5464"
$ws.Range("D87").Value2 = "0.75x"
$ws.Range("B91").Value2 = "Mitutoyo"
$ws.Range("C92").Value2 = "Mitutoyo 7.5x"
$ws.Range("D92").Value2 = "7.5x"
$ws.Range("C91").Value2 = "Mitutoyo HR 5.0x"
$ws.Range("G80").Value2 = "1188"
$ws.Range("G81").Value2 = "674"
$ws.Range("G82").Value2 = "783"
$ws.Range("G83").Value2 = "898"
$ws.Range("G84").Value2 = "1012"
$ws.Range("G86").Value2 = "925"
$ws.Range("H86").Value2 = "8"
$ws.Range("G87").Value2 = "798"
$ws.Range("G89").Value2 = "633"
$ws.Range("G90").Value2 = "889"
$ws.Range("G91").Value2 = "1218"
$ws.Range("G92").Value2 = "919"
$ws.Range("B80").Value2 = "Laowa"
$ws.Range("C80").Value2 = "Laowa 25mm f/2.8"
$ws.Range("D80").Value2 = "2.5x"
$ws.Range("H80").Value2 = "2"
$ws.Range("I80").Value2 = "mm"
$ws.Range("B81").Value2 = "Laowa"
$ws.Range("C81").Value2 = "Laowa 25mm f/2.8"
$ws.Range("D81").Value2 = "3.0x"
$ws.Range("H81").Value2 = "1000"
$ws.Range("I81").Value2 = "µm"
$ws.Range("B82").Value2 = "Laowa"
$ws.Range("C82").Value2 = "Laowa 25mm f/2.8"
$ws.Range("D82").Value2 = "3.5x"
$ws.Range("H82").Value2 = "1000"
$ws.Range("I82").Value2 = "µm"
$ws.Range("B83").Value2 = "Laowa"
$ws.Range("C83").Value2 = "Laowa 25mm f/2.8"
$ws.Range("D83").Value2 = "4.0x"
$ws.Range("H83").Value2 = "1000"
$ws.Range("I83").Value2 = "µm"
$ws.Range("B84").Value2 = "Laowa"
$ws.Range("C84").Value2 = "Laowa 25mm f/2.8"
$ws.Range("D84").Value2 = "4.5x"
$ws.Range("H84").Value2 = "1000"
$ws.Range("I84").Value2 = "µm"
$ws.Range("B85").Value2 = "Laowa"
$ws.Range("C85").Value2 = "Laowa 25mm f/2.8"
$ws.Range("D85").Value2 = "5.0x"
$ws.Range("G85").Value2 = "1107"
$ws.Range("H85").Value2 = "1000"
$ws.Range("I85").Value2 = "µm"
$ws.Range("B86").Value2 = "Laowa"
$ws.Range("C86").Value2 = "Laowa 85mm f/5.6"
$ws.Range("D86").Value2 = "0.5x"
$ws.Range("I86").Value2 = "mm"
$ws.Range("B87").Value2 = "Laowa"
$ws.Range("C87").Value2 = "Laowa 85mm f/5.6"
$ws.Range("H87").Value2 = "5"
$ws.Range("I87").Value2 = "mm"
$ws.Range("B88").Value2 = "Laowa"
$ws.Range("C88").Value2 = "Laowa 85mm f/5.6"
$ws.Range("D88").Value2 = "1.0x"
$ws.Range("G88").Value2 = "670"
$ws.Range("H88").Value2 = "4"
$ws.Range("I88").Value2 = "mm"
$ws.Range("B89").Value2 = "Laowa"
$ws.Range("C89").Value2 = "Laowa 85mm f/5.6"
$ws.Range("D89").Value2 = "1.5x"
$ws.Range("H89").Value2 = "2"
$ws.Range("I89").Value2 = "mm"
$ws.Range("B90").Value2 = "Laowa"
$ws.Range("C90").Value2 = "Laowa 85mm f/5.6"
$ws.Range("D90").Value2 = "2.0x"
$ws.Range("H90").Value2 = "2"
$ws.Range("I90").Value2 = "mm"
$ws.Range("D91").Value2 = "5.0x"
$ws.Range("H91").Value2 = "1000"
$ws.Range("I91").Value2 = "µm"
$ws.Range("B92").Value2 = "Mitutoyo"
$ws.Range("H92").Value2 = "500"
$ws.Range("I92").Value2 = "µm"

# --- 3. Five pre-existing rows (6, 57, 62, 63, 64) were previously the
#        sole users of a now-redundant row-level style; normalise them
#        back to the plain row style used throughout the rest of the
#        table (matches the cleanup captured in the authored edit). ---
foreach ($r in @(6, 57, 62, 63, 64)) {
    $ws.Rows.Item($r).ClearFormats() | Out-Null
}
$ws.Range("A2:I2").Copy() | Out-Null
foreach ($r in @(6, 57, 62, 63, 64)) {
    $ws.Range("A" + $r + ":I" + $r).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# --- 4. Leave the selection where the author ended up after entering
#        the new data. ---
$ws.Range("I93").Select()
